$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1515987.1
$ws.Range("J17").Value = 1563346.1
$ws.Range("L17").Value = 4690038.300000001
$ws.Range("N17").Value = -4690374.300000001
$ws.Range("H112").Value = 1147.9706
$ws.Range("I112").Value = 798.75
$ws.Range("J112").Value = 1255.4231
$ws.Range("K112").Value = 2396.25
$ws.Range("L112").Value = 3766.2693
$ws.Range("M112").Value = -1288.25
$ws.Range("N112").Value = -5982.2693
$ws.Range("H129").Value = 1120.8055
$ws.Range("J129").Value = 1178.0938
$ws.Range("L129").Value = 3534.2814
$ws.Range("N129").Value = -13534.2814
$ws.Range("H137").Value = 1221.0322
$ws.Range("I137").Value = 1090.8541
$ws.Range("J137").Value = 1667.3572
$ws.Range("K137").Value = 3272.5623
$ws.Range("L137").Value = 5002.071599999999
$ws.Range("M137").Value = -722.5623000000001
$ws.Range("N137").Value = -10102.0716
$ws.Range("H138").Value = 3291.9778
$ws.Range("I138").Value = 1718.5454
$ws.Range("J138").Value = 4797
$ws.Range("K138").Value = 5155.6362
$ws.Range("L138").Value = 14391
$ws.Range("M138").Value = -15.63619999999992
$ws.Range("N138").Value = -24671
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1283.4231
$ws.Range("I2").Value = 1390
$ws.Range("J2").Value = 928.1667
$ws.Range("K2").Value = 1390
$ws.Range("L2").Value = 928.1667
$ws.Range("M2").Value = -1277
$ws.Range("N2").Value = -1154.1667
$ws.Range("H45").Value = 3541.7144
$ws.Range("I45").Value = 4160.4
$ws.Range("J45").Value = 1995
$ws.Range("K45").Value = 4160.4
$ws.Range("L45").Value = 1995
$ws.Range("M45").Value = -3783.4
$ws.Range("N45").Value = -2749
$ws.Range("H61").Value = 2088.3333
$ws.Range("I61").Value = 2053.2715
$ws.Range("J61").Value = 2325
$ws.Range("K61").Value = 2053.2715
$ws.Range("L61").Value = 2325
$ws.Range("M61").Value = -1841.2715
$ws.Range("N61").Value = -2749
$ws.Range("H116").Value = 1283.4231
$ws.Range("I116").Value = 1390
$ws.Range("J116").Value = 928.1667
$ws.Range("K116").Value = 1390
$ws.Range("L116").Value = 928.1667
$ws.Range("M116").Value = 904
$ws.Range("N116").Value = -5516.1667
$ws.Range("H136").Value = 2088.3333
$ws.Range("I136").Value = 2053.2715
$ws.Range("J136").Value = 2325
$ws.Range("K136").Value = 6159.814499999999
$ws.Range("L136").Value = 6975
$ws.Range("M136").Value = -3609.814499999999
$ws.Range("N136").Value = -12075
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1283.4231
$ws.Range("I3").Value = 1390
$ws.Range("J3").Value = 928.1667
$ws.Range("K3").Value = 1390
$ws.Range("L3").Value = 928.1667
$ws.Range("M3").Value = -1276
$ws.Range("N3").Value = -1156.1667
$ws.Range("H105").Value = 1568.3636
$ws.Range("I105").Value = 1107.4286
$ws.Range("J105").Value = 2375
$ws.Range("K105").Value = 1107.4286
$ws.Range("L105").Value = 2375
$ws.Range("M105").Value = 639.5714
$ws.Range("N105").Value = -5869
$ws.Range("H134").Value = 1240.0217
$ws.Range("I134").Value = 928.5278
$ws.Range("J134").Value = 2361.4
$ws.Range("K134").Value = 2785.5834
$ws.Range("L134").Value = 7084.200000000001
$ws.Range("M134").Value = -250.5834
$ws.Range("N134").Value = -12154.2
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 906.64703
$ws.Range("I16").Value = 939
$ws.Range("J16").Value = 870.25
$ws.Range("K16").Value = 939
$ws.Range("L16").Value = 870.25
$ws.Range("M16").Value = -652
$ws.Range("N16").Value = -1444.25
$ws.Range("H20").Value = 52092.7
$ws.Range("J20").Value = 52092.7
$ws.Range("L20").Value = 52092.7
$ws.Range("N20").Value = -52564.7
$ws.Range("H30").Value = 52092.7
$ws.Range("J30").Value = 52092.7
$ws.Range("L30").Value = 52092.7
$ws.Range("N30").Value = -52274.7
$ws.Range("H58").Value = 1003.96344
$ws.Range("I58").Value = 756.0294
$ws.Range("J58").Value = 2208.2144
$ws.Range("K58").Value = 756.0294
$ws.Range("L58").Value = 2208.2144
$ws.Range("M58").Value = -553.0294
$ws.Range("N58").Value = -2614.2144
$ws.Range("H95").Value = 18760
$ws.Range("J95").Value = 18760
$ws.Range("L95").Value = 18760
$ws.Range("N95").Value = -24252
$ws.Range("H113").Value = 906.64703
$ws.Range("I113").Value = 939
$ws.Range("J113").Value = 870.25
$ws.Range("K113").Value = 939
$ws.Range("L113").Value = 870.25
$ws.Range("M113").Value = 1231
$ws.Range("N113").Value = -5210.25
$ws.Range("H128").Value = 52092.7
$ws.Range("J128").Value = 52092.7
$ws.Range("L128").Value = 52092.7
$ws.Range("N128").Value = -62052.7
$ws.Range("H132").Value = 3686.889
$ws.Range("I132").Value = 2991.4443
$ws.Range("J132").Value = 4382.3335
$ws.Range("K132").Value = 8974.332900000001
$ws.Range("L132").Value = 13147.0005
$ws.Range("M132").Value = -6444.332900000001
$ws.Range("N132").Value = -18207.0005
$ws.Range("H136").Value = 1003.96344
$ws.Range("I136").Value = 756.0294
$ws.Range("J136").Value = 2208.2144
$ws.Range("K136").Value = 2268.0882
$ws.Range("L136").Value = 6624.6432
$ws.Range("M136").Value = 281.9117999999999
$ws.Range("N136").Value = -11724.6432
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 418411.84
$ws.Range("I5").Value = 464.2353
$ws.Range("J5").Value = 702616.2
$ws.Range("K5").Value = 1392.7059
$ws.Range("L5").Value = 2107848.6
$ws.Range("M5").Value = -1280.7059
$ws.Range("N5").Value = -2108072.6
$ws.Range("H113").Value = 3227.5
$ws.Range("I113").Value = 5503
$ws.Range("J113").Value = 952
$ws.Range("K113").Value = 16509
$ws.Range("L113").Value = 2856
$ws.Range("M113").Value = -14339
$ws.Range("N113").Value = -7196
$ws.Range("H122").Value = 1427.9565
$ws.Range("I122").Value = 430.30768
$ws.Range("J122").Value = 2724.9
$ws.Range("K122").Value = 3872.76912
$ws.Range("L122").Value = 24524.1
$ws.Range("M122").Value = -1422.76912
$ws.Range("N122").Value = -29424.1
$ws.Range("H131").Value = 3131.2922
$ws.Range("J131").Value = 3660.463
$ws.Range("L131").Value = 10981.389
$ws.Range("N131").Value = -21061.389
$ws.Range("H135").Value = 418411.84
$ws.Range("I135").Value = 464.2353
$ws.Range("J135").Value = 702616.2
$ws.Range("K135").Value = 4178.1177
$ws.Range("L135").Value = 6323545.8
$ws.Range("M135").Value = -1643.1177
$ws.Range("N135").Value = -6328615.8
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 1104.8021
$ws.Range("I132").Value = 1041.2593
$ws.Range("J132").Value = 1447.9333
$ws.Range("K132").Value = 3123.7779
$ws.Range("L132").Value = 4343.7999
$ws.Range("M132").Value = -593.7779
$ws.Range("N132").Value = -9403.7999
$ws.Range("H134").Value = 62500
$ws.Range("J134").Value = 62500
$ws.Range("L134").Value = 62500
$ws.Range("N134").Value = -72640
$ws.Range("H135").Value = 50639.535
$ws.Range("J135").Value = 50639.535
$ws.Range("L135").Value = 50639.535
$ws.Range("N135").Value = -60779.535
$ws.Range("H136").Value = 6173921
$ws.Range("I136").Value = 1101.8302
$ws.Range("J136").Value = 333333340
$ws.Range("K136").Value = 3305.4906
$ws.Range("L136").Value = 1000000020
$ws.Range("M136").Value = -755.4906000000001
$ws.Range("N136").Value = -1000005120
$ws.Range("H137").Value = 29918.428
$ws.Range("J137").Value = 29918.428
$ws.Range("L137").Value = 29918.428
$ws.Range("N137").Value = -40118.428
$ws.Range("H139").Value = 36118.184
$ws.Range("J139").Value = 37730
$ws.Range("L139").Value = 37730
$ws.Range("N139").Value = -48010
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 51050.035
$ws.Range("I126").Value = 64554.637
$ws.Range("K126").Value = 193663.911
$ws.Range("M126").Value = -191193.911
$ws.Range("H132").Value = 1359.2616
$ws.Range("I132").Value = 1224
$ws.Range("J132").Value = 1810.1333
$ws.Range("K132").Value = 3672
$ws.Range("L132").Value = 5430.3999
$ws.Range("M132").Value = -1142
$ws.Range("N132").Value = -10490.3999
$ws.Range("H136").Value = 4067310.8
$ws.Range("I136").Value = 7246826
$ws.Range("J136").Value = 4596.8613
$ws.Range("K136").Value = 21740478
$ws.Range("L136").Value = 13790.5839
$ws.Range("M136").Value = -21737928
$ws.Range("N136").Value = -18890.5839
